# Update gh-pages to output generated at 456a3b4
# Bumps the "想去人数" (want-to-go count) figures in column F across the
# 展览 (Exhibitions), 演出 (Performances) and 全部类型 (All types) sheets.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet -------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value  = 3
$wsExpo.Range("F6").Value  = 193
$wsExpo.Range("F7").Value  = 4071
$wsExpo.Range("F9").Value  = 118
$wsExpo.Range("F11").Value = 86
$wsExpo.Range("F12").Value = 85
$wsExpo.Range("F13").Value = 685
$wsExpo.Range("F15").Value = 955
$wsExpo.Range("F16").Value = 72
$wsExpo.Range("F17").Value = 236
$wsExpo.Range("F20").Value = 107
$wsExpo.Range("F21").Value = 92
$wsExpo.Range("F22").Value = 3434
$wsExpo.Range("F23").Value = 5769
$wsExpo.Range("F25").Value = 28
$wsExpo.Range("F26").Value = 83
$wsExpo.Range("F27").Value = 518
$wsExpo.Range("F29").Value = 3342
$wsExpo.Range("F30").Value = 354
$wsExpo.Range("F32").Value = 2449
$wsExpo.Range("F35").Value = 119
$wsExpo.Range("F36").Value = 202
$wsExpo.Range("F37").Value = 257
$wsExpo.Range("F38").Value = 345
$wsExpo.Range("F40").Value = 1006
$wsExpo.Range("F42").Value = 9

# --- 演出 (Performances) sheet ------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 94

# --- 全部类型 (All types) sheet -----------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value  = 3
$wsAll.Range("F6").Value  = 193
$wsAll.Range("F7").Value  = 4072
$wsAll.Range("F9").Value  = 118
$wsAll.Range("F11").Value = 94
$wsAll.Range("F12").Value = 86
$wsAll.Range("F13").Value = 85
$wsAll.Range("F14").Value = 685
$wsAll.Range("F16").Value = 955
$wsAll.Range("F17").Value = 72
$wsAll.Range("F18").Value = 236
$wsAll.Range("F21").Value = 107
$wsAll.Range("F22").Value = 92
$wsAll.Range("F23").Value = 3434
$wsAll.Range("F24").Value = 5769
$wsAll.Range("F26").Value = 28
$wsAll.Range("F27").Value = 83
$wsAll.Range("F28").Value = 518
$wsAll.Range("F30").Value = 3342
$wsAll.Range("F31").Value = 354
$wsAll.Range("F33").Value = 2449
$wsAll.Range("F36").Value = 119
$wsAll.Range("F37").Value = 202
$wsAll.Range("F38").Value = 257
$wsAll.Range("F39").Value = 345
$wsAll.Range("F41").Value = 1006
$wsAll.Range("F43").Value = 9
